# Update the table-of-contents page numbers (column D) to reflect the
# final pagination of the thesis, then leave the selection on the page
# number column (K8:K28) as the author did before sending.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pageNumbers = @{
    8  = "5"
    9  = "6"
    10 = "6"
    11 = "10"
    12 = "12"
    13 = "12"
    14 = "14"
    15 = "22"
    16 = "22"
    17 = "23"
    18 = "33"
    19 = "36"
    20 = "36"
    21 = "37"
    22 = "43"
    23 = "52"
    24 = "53"
    25 = "53"
    26 = "54"
    27 = "59"
    28 = "60"
}

foreach ($row in $pageNumbers.Keys) {
    $ws.Range("D$row").Value = $pageNumbers[$row]
}

# K8 used to hold its own (non-shared) formula "=D8"; it now holds the
# resolved page number as a literal value, same as the other K cells'
# cached results.
$ws.Range("K8").Value = $pageNumbers[8]

$ws.Range("K8:K28").Select()
